$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the metrics_file value (F2) first so the new shared string is registered
# in the same order as the saved workbook.
$ws.Range("F2").Value = "datasets/10x_visiumhd_mouse_brain/metrics_summary.csv"

# Update path column (E2) to point at the renamed dataset folder and matrix file
$ws.Range("E2").Value = "datasets/10x_visiumhd_mouse_brain/segmented_outputs/filtered_feature_cell_matrix.h5"

# Update the active selection to match the saved view
$ws.Range("E3").Select()
